$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3F")

# Fix C3: was stored as inline text "28" -> should become a real number
$ws.Range("C3").Value = 28

# New submission row synced: 2026-02-08 22:15:10
$ws.Range("A4").Value = "2026-02-08 22:15:10"
$ws.Range("B4").Value = "Usman Muhammad Gubio"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "05"
$ws.Range("D4").Value = 7
